$wb = $excel.ActiveWorkbook

# --- Primers_for_uha ---
$ws = $wb.Worksheets.Item('Primers_for_uha')
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("B3").Value = 'aceE_del'
$ws.Range("C3").Value = 'GCGTCACAGACATGAAATTGGT'
$ws.Range("D3").Value = 'AGCCATTATTCTTTTACCTCGGGTTATTCCTTATCTATCT'
$ws.Range("E3").Value = 'GCGTCACAGACATGAAATTGGTAAGACCAATTGACTTCGGCAAGTGGCTTAAGACAGGAACTCATGGCCTACAGCAAAATCCGCCAACCAAAACTCTCCGATGTGATTGAGCAGCAACTGGAGTTTTTGATCCTCGAAGGCACTCTCCGCCCGGGCGAAAAACTCCCACCGGAACGCGAACTGGCAAAACAGTTTGACGTCTCCCGTCCCTCCTTGCGTGAGGCGATTCAACGTCTCGAAGCGAAGGGCTTGTTGCTTCGTCGCCAGGGTGGCGGCACTTTTGTCCAGAGCAGCCTATGGCAAAGCTTCAGCGATCCGCTGGTGGAGCTGCTCTCCGACCATCCTGAGTCACAGTATGACTTGCTCGAAACACGACACGCCCTGGAAGGTATCGCCGCTTATTACGCCGCGCTGCGTAGTACCGATGAAGACAAGGAACGCATCCGTGAACTCCACCACGCCATAGAGCTGGCGCAGCAGTCTGGCGATCTGGACGCGGAATCAAACGCCGTACTCCAGTATCAGATTGCCGTCACCGAAGCGGCCCACAATGTGGTTCTGCTTCATCTGCTAAGGTGTATGGAGCCGATGTTGGCCCAGAATGTCCGCCAGAACTTCGAATTGCTCTATTCGCGTCGCGAGATGCTGCCGCTGGTGAGTAGTCACCGCACCCGCATATTTGAAGCGATTATGGCCGGTAAGCCGGAAGAAGCGCGCGAAGCATCGCATCGCCATCTGGCCTTTATCGAAGAAATTTTGCTCGACAGAAGTCGTGAAGAGAGCCGCCGTGAGCGTTCTCTGCGTCGTCTGGAGCAACGAAAGAATTAGTGATTTTTCTGGTAAAAATTATCCAGAAGATGTTGTAAATCAAGCGCATATAAAAGCGCGGCAACTAAACGTAGAACCTGTCTTATTGAGCTTTCCGGCGAGAGTTCAATGGGACAGGTTCCAGAAAACTCAACGTTATTAGATAGATAAGGAATAACCC'
$ws.Range("F3").Value = 'TACGTAAAGTCTACATTTGTGCATAGTTACAACTTTGAAACGTTATATATGTCAAGTTGTTAAAATGTGCACAGTTTCATGATTTCAATCAAAACCTGTATGGACATAAGGTGAATACTTTGTTACTTTAGCGTCACAGACATGAAATTGGTAAGACCAATTGACTTCGGCAAGTGGCTTAAGACAGGAACTCATGGCCTACAGCAAAATCCGCCAACCAAAACTCTCCGATGTGATTGAGCAGCAACTGGAGTTTTTGATCCTCGAAGGCACTCTCCGCCCGGGCGAAAAACTCCCACCGGAACGCGAACTGGCAAAACAGTTTGACGTCTCCCGTCCCTCCTTGCGTGAGGCGATTCAACGTCTCGAAGCGAAGGGCTTGTTGCTTCGTCGCCAGGGTGGCGGCACTTTTGTCCAGAGCAGCCTATGGCAAAGCTTCAGCGATCCGCTGGTGGAGCTGCTCTCCGACCATCCTGAGTCACAGTATGACTTGCTCGAAACACGACACGCCCTGGAAGGTATCGCCGCTTATTACGCCGCGCTGCGTAGTACCGATGAAGACAAGGAACGCATCCGTGAACTCCACCACGCCATAGAGCTGGCGCAGCAGTCTGGCGATCTGGACGCGGAATCAAACGCCGTACTCCAGTATCAGATTGCCGTCACCGAAGCGGCCCACAATGTGGTTCTGCTTCATCTGCTAAGGTGTATGGAGCCGATGTTGGCCCAGAATGTCCGCCAGAACTTCGAATTGCTCTATTCGCGTCGCGAGATGCTGCCGCTGGTGAGTAGTCACCGCACCCGCATATTTGAAGCGATTATGGCCGGTAAGCCGGAAGAAGCGCGCGAAGCATCGCATCGCCATCTGGCCTTTATCGAAGAAATTTTGCTCGACAGAAGTCGTGAAGAGAGCCGCCGTGAGCGTTCTCTGCGTCGTCTGGAGCAACGAAAGAATTAGTGATTTTTCTGGTAAAAATTATCCAGAAGATGTTGTAAATCAAGCGCATATAAAAGCGCGGCAACTAAACGTAGAACCTGTCTTATTGAGCTTTCCGGCGAGAGTTCAATGGGACAGGTTCCAGAAAACTCAACGTTATTAGATAGATAAGGAATAACCC'
$ws.Range("G3").Value = 'GCGTCACAGACATGAAATTGGTAAGACCAATTGACTTCGGCAAGTGGCTTAAGACAGGAACTCATGGCCTACAGCAAAATCCGCCAACCAAAACTCTCCGATGTGATTGAGCAGCAACTGGAGTTTTTGATCCTCGAAGGCACTCTCCGCCCGGGCGAAAAACTCCCACCGGAACGCGAACTGGCAAAACAGTTTGACGTCTCCCGTCCCTCCTTGCGTGAGGCGATTCAACGTCTCGAAGCGAAGGGCTTGTTGCTTCGTCGCCAGGGTGGCGGCACTTTTGTCCAGAGCAGCCTATGGCAAAGCTTCAGCGATCCGCTGGTGGAGCTGCTCTCCGACCATCCTGAGTCACAGTATGACTTGCTCGAAACACGACACGCCCTGGAAGGTATCGCCGCTTATTACGCCGCGCTGCGTAGTACCGATGAAGACAAGGAACGCATCCGTGAACTCCACCACGCCATAGAGCTGGCGCAGCAGTCTGGCGATCTGGACGCGGAATCAAACGCCGTACTCCAGTATCAGATTGCCGTCACCGAAGCGGCCCACAATGTGGTTCTGCTTCATCTGCTAAGGTGTATGGAGCCGATGTTGGCCCAGAATGTCCGCCAGAACTTCGAATTGCTCTATTCGCGTCGCGAGATGCTGCCGCTGGTGAGTAGTCACCGCACCCGCATATTTGAAGCGATTATGGCCGGTAAGCCGGAAGAAGCGCGCGAAGCATCGCATCGCCATCTGGCCTTTATCGAAGAAATTTTGCTCGACAGAAGTCGTGAAGAGAGCCGCCGTGAGCGTTCTCTGCGTCGTCTGGAGCAACGAAAGAATTAGTGATTTTTCTGGTAAAAATTATCCAGAAGATGTTGTAAATCAAGCGCATATAAAAGCGCGGCAACTAAACGTAGAACCTGTCTTATTGAGCTTTCCGGCGAGAGTTCAATGGGACAGGTTCCAGAAAACTCAACGTTATTAGATAGATAAGGAATAACCCGAGGTAAAAGAATAATGGCT'
$ws.Range("H3").Value = 1008

# --- Primers_for_dha ---
$ws = $wb.Worksheets.Item('Primers_for_dha')
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("B3").Value = 'aceE_del'
$ws.Range("C3").Value = 'GAGGTAAAAGAATAATGGCTATCGA'
$ws.Range("D3").Value = 'TTAACACCAAACTCGCGTGC'
$ws.Range("E3").Value = 'GAGGTAAAAGAATAATGGCTATCGAAATCAAAGTACCGGACATCGGGGCTGATGAAGTTGAAATCACCGAGATCCTGGTCAAAGTGGGCGACAAAGTTGAAGCCGAACAGTCGCTGATCACCGTAGAAGGCGACAAAGCCTCTATGGAAGTTCCGTCTCCGCAGGCGGGTATCGTTAAAGAGATCAAAGTCTCTGTTGGCGATAAAACCCAGACCGGCGCACTGATTATGATTTTCGATTCCGCCGACGGTGCAGCAGACGCTGCACCTGCTCAGGCAGAAGAGAAGAAAGAAGCAGCTCCGGCAGCAGCACCAGCGGCTGCGGCGGCAAAAGACGTTAACGTTCCGGATATCGGCAGCGACGAAGTTGAAGTGACCGAAATCCTGGTGAAAGTTGGCGATAAAGTTGAAGCTGAACAGTCGCTGATCACCGTAGAAGGCGACAAGGCTTCTATGGAAGTTCCGGCTCCGTTTGCTGGCACCGTGAAAGAGATCAAAGTGAACGTGGGTGACAAAGTGTCTACCGGCTCGCTGATTATGGTCTTCGAAGTCGCGGGTGAAGCAGGCGCGGCAGCTCCGGCCGCTAAACAGGAAGCAGCTCCGGCAGCGGCCCCTGCACCAGCGGCTGGCGTGAAAGAAGTTAACGTTCCGGATATCGGCGGTGACGAAGTTGAAGTGACTGAAGTGATGGTGAAAGTGGGCGACAAAGTTGCCGCTGAACAGTCACTGATCACCGTAGAAGGCGACAAAGCTTCTATGGAAGTTCCGGCGCCGTTTGCAGGCGTCGTGAAGGAACTGAAAGTCAACGTTGGCGATAAAGTGAAAACTGGCTCGCTGATTATGATCTTCGAAGTTGAAGGCGCAGCGCCTGCGGCAGCTCCTGCGAAACAGGAAGCGGCAGCGCCGGCACCGGCAGCAAAAGCTGAAGCCCCGGCAGCAGCACCAGCTGCGAAAGCGGAAGGCAAATCTGAATTTGCTGAAAACGACGCTTATGTTCACGCGACTCCGCTGATCCGCCGTCTGGCACGCGAGTTTGGTGTTAA'
$ws.Range("F3").Value = 1042
$ws.Range("G3").Value = 'GAGGTAAAAGAATAATGGCTATCGAAATCAAAGTACCGGACATCGGGGCTGATGAAGTTGAAATCACCGAGATCCTGGTCAAAGTGGGCGACAAAGTTGAAGCCGAACAGTCGCTGATCACCGTAGAAGGCGACAAAGCCTCTATGGAAGTTCCGTCTCCGCAGGCGGGTATCGTTAAAGAGATCAAAGTCTCTGTTGGCGATAAAACCCAGACCGGCGCACTGATTATGATTTTCGATTCCGCCGACGGTGCAGCAGACGCTGCACCTGCTCAGGCAGAAGAGAAGAAAGAAGCAGCTCCGGCAGCAGCACCAGCGGCTGCGGCGGCAAAAGACGTTAACGTTCCGGATATCGGCAGCGACGAAGTTGAAGTGACCGAAATCCTGGTGAAAGTTGGCGATAAAGTTGAAGCTGAACAGTCGCTGATCACCGTAGAAGGCGACAAGGCTTCTATGGAAGTTCCGGCTCCGTTTGCTGGCACCGTGAAAGAGATCAAAGTGAACGTGGGTGACAAAGTGTCTACCGGCTCGCTGATTATGGTCTTCGAAGTCGCGGGTGAAGCAGGCGCGGCAGCTCCGGCCGCTAAACAGGAAGCAGCTCCGGCAGCGGCCCCTGCACCAGCGGCTGGCGTGAAAGAAGTTAACGTTCCGGATATCGGCGGTGACGAAGTTGAAGTGACTGAAGTGATGGTGAAAGTGGGCGACAAAGTTGCCGCTGAACAGTCACTGATCACCGTAGAAGGCGACAAAGCTTCTATGGAAGTTCCGGCGCCGTTTGCAGGCGTCGTGAAGGAACTGAAAGTCAACGTTGGCGATAAAGTGAAAACTGGCTCGCTGATTATGATCTTCGAAGTTGAAGGCGCAGCGCCTGCGGCAGCTCCTGCGAAACAGGAAGCGGCAGCGCCGGCACCGGCAGCAAAAGCTGAAGCCCCGGCAGCAGCACCAGCTGCGAAAGCGGAAGGCAAATCTGAATTTGCTGAAAACGACGCTTATGTTCACGCGACTCCGCTGATCCGCCGTCTGGCACGCGAGTTTGGTGTTAACCTTGCGAAAGTGAAGGGCACTGGCCGTAAAGGTCGTATCCTGCGCGAAGACGTTCAGGCTTACGTGAAAGAAGCTATCAAACGTGCAGAAGCAGCTCCGGCAGCGACTGGCGGTGGTATCCCTGGCATG'
$ws.Range("H3").Value = 'GAGGTAAAAGAATAATGGCTATCGAAATCAAAGTACCGGACATCGGGGCTGATGAAGTTGAAATCACCGAGATCCTGGTCAAAGTGGGCGACAAAGTTGAAGCCGAACAGTCGCTGATCACCGTAGAAGGCGACAAAGCCTCTATGGAAGTTCCGTCTCCGCAGGCGGGTATCGTTAAAGAGATCAAAGTCTCTGTTGGCGATAAAACCCAGACCGGCGCACTGATTATGATTTTCGATTCCGCCGACGGTGCAGCAGACGCTGCACCTGCTCAGGCAGAAGAGAAGAAAGAAGCAGCTCCGGCAGCAGCACCAGCGGCTGCGGCGGCAAAAGACGTTAACGTTCCGGATATCGGCAGCGACGAAGTTGAAGTGACCGAAATCCTGGTGAAAGTTGGCGATAAAGTTGAAGCTGAACAGTCGCTGATCACCGTAGAAGGCGACAAGGCTTCTATGGAAGTTCCGGCTCCGTTTGCTGGCACCGTGAAAGAGATCAAAGTGAACGTGGGTGACAAAGTGTCTACCGGCTCGCTGATTATGGTCTTCGAAGTCGCGGGTGAAGCAGGCGCGGCAGCTCCGGCCGCTAAACAGGAAGCAGCTCCGGCAGCGGCCCCTGCACCAGCGGCTGGCGTGAAAGAAGTTAACGTTCCGGATATCGGCGGTGACGAAGTTGAAGTGACTGAAGTGATGGTGAAAGTGGGCGACAAAGTTGCCGCTGAACAGTCACTGATCACCGTAGAAGGCGACAAAGCTTCTATGGAAGTTCCGGCGCCGTTTGCAGGCGTCGTGAAGGAACTGAAAGTCAACGTTGGCGATAAAGTGAAAACTGGCTCGCTGATTATGATCTTCGAAGTTGAAGGCGCAGCGCCTGCGGCAGCTCCTGCGAAACAGGAAGCGGCAGCGCCGGCACCGGCAGCAAAAGCTGAAGCCCCGGCAGCAGCACCAGCTGCGAAAGCGGAAGGCAAATCTGAATTTGCTGAAAACGACGCTTATGTTCACGCGACTCCGCTGATCCGCCGTCTGGCACGCGAGTTTGGTGTTAA'

# --- Primers_for_verify2 ---
$ws = $wb.Worksheets.Item('Primers_for_verify2')
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("B3").Value = 'aceE_del'
$ws.Range("C3").Value = 'ACGTAAAGTCTACATTTGTGCATAGTTACAACTTTGAAACGTTATATATGTCAAGTTGTTAAAATGTGCACAGTTTCATGATTTCAATCAAAACCTGTATGGACATAAGGTGAATACTTTGTTACTTTAGCGTCACAGACATGAAATTGGTAAGACCAATTGACTTCGGCAAGTGGCTTAAGACAGGAACTCATGGCCTACAGCAAAATCCGCCAACCAAAACTCTCCGATGTGATTGAGCAGCAACTGGAGTTTTTGATCCTCGAAGGCACTCTCCGCCCGGGCGAAAAACTCCCACCGGAACGCGAACTGGCAAAACAGTTTGACGTCTCCCGTCCCTCCTTGCGTGAGGCGATTCAACGTCTCGAAGCGAAGGGCTTGTTGCTTCGTCGCCAGGGTGGCGGCACTTTTGTCCAGAGCAGCCTATGGCAAAGCTTCAGCGATCCGCTGGTGGAGCTGCTCTCCGACCATCCTGAGTCACAGTATGACTTGCTCGAAACACGACACGCCCTGGAAGGTATCGCCGCTTATTACGCCGCGCTGCGTAGTACCGATGAAGACAAGGAACGCATCCGTGAACTCCACCACGCCATAGAGCTGGCGCAGCAGTCTGGCGATCTGGACGCGGAATCAAACGCCGTACTCCAGTATCAGATTGCCGTCACCGAAGCGGCCCACAATGTGGTTCTGCTTCATCTGCTAAGGTGTATGGAGCCGATGTTGGCCCAGAATGTCCGCCAGAACTTCGAATTGCTCTATTCGCGTCGCGAGATGCTGCCGCTGGTGAGTAGTCACCGCACCCGCATATTTGAAGCGATTATGGCCGGTAAGCCGGAAGAAGCGCGCGAAGCATCGCATCGCCATCTGGCCTTTATCGAAGAAATTTTGCTCGACAGAAGTCGTGAAGAGAGCCGCCGTGAGCGTTCTCTGCGTCGTCTGGAGCAACGAAAGAATTAGTGATTTTTCTGGTAAAAATTATCCAGAAGATGTTGTAAATCAAGCGCATATAAAAGCGCGGCAACTAAACGTAGAACCTGTCTTATTGAGCTTTCCGGCGAGAGTTCAATGGGACAGGTTCCAGAAAACTCAACGTTATTAGATAGATAAGGAATAACCCGAGGTAAAAGAATAATGGCTATCGAAATCAAAGTACCGGACATCGGGGCTGATGAAGTTGAAATCACCGAGATCCTGGTCAAAGTGGGCGACAAAGTTGAAGCCGAACAGTCGCTGATCACCGTAGAAGGCGACAAAGCCTCTATGGAAGTTCCGTCTCCGCAGGCGGGTATCGTTAAAGAGATCAAAGTCTCTGTTGGCGATAAAACCCAGACCGGCGCACTGATTATGATTTTCGATTCCGCCGACGGTGCAGCAGACGCTGCACCTGCTCAGGCAGAAGAGAAGAAAGAAGCAGCTCCGGCAGCAGCACCAGCGGCTGCGGCGGCAAAAGACGTTAACGTTCCGGATATCGGCAGCGACGAAGTTGAAGTGACCGAAATCCTGGTGAAAGTTGGCGATAAAGTTGAAGCTGAACAGTCGCTGATCACCGTAGAAGGCGACAAGGCTTCTATGGAAGTTCCGGCTCCGTTTGCTGGCACCGTGAAAGAGATCAAAGTGAACGTGGGTGACAAAGTGTCTACCGGCTCGCTGATTATGGTCTTCGAAGTCGCGGGTGAAGCAGGCGCGGCAGCTCCGGCCGCTAAACAGGAAGCAGCTCCGGCAGCGGCCCCTGCACCAGCGGCTGGCGTGAAAGAAGTTAACGTTCCGGATATCGGCGGTGACGAAGTTGAAGTGACTGAAGTGATGGTGAAAGTGGGCGACAAAGTTGCCGCTGAACAGTCACTGATCACCGTAGAAGGCGACAAAGCTTCTATGGAAGTTCCGGCGCCGTTTGCAGGCGTCGTGAAGGAACTGAAAGTCAACGTTGGCGATAAAGTGAAAACTGGCTCGCTGATTATGATCTTCGAAGTTGAAGGCGCAGCGCCTGCGGCAGCTCCTGCGAAACAGGAAGCGGCAGCGCCGGCACCGGCAGCAAAAGCTGAAGCCCCGGCAGCAGCACCAGCTGCGAAAGCGGAAGGCAAATCTGAATTTGCTGAAAACGACGCTTATGTTCACGCGACTCCGCTGATCCGCCGTCTGGCACGCGAGTTTGGTGTTAACCTTGCGAAAGTGAAGGGCACTGGCCGTAAAGGTCGTATCCTGCGCGAAGACGTTCAGGCTTACGTGAAAGAAGCTATCAAACGTGCAGAAGCAGCTCC'
$ws.Range("D3").Value = -42.60669801557233
$ws.Range("E3").Value = -0.900577242230667
$ws.Range("F3").Value = -6.517879085020752
$ws.Range("G3").Value = 'ACGTAAAGTCTACATTTGTGCA'
$ws.Range("H3").Value = 'GGAGCTGCTTCTGCACGTTT'
$ws.Range("I3").Value = 56.32043528895827
$ws.Range("J3").Value = 61.51048014864102
$ws.Range("K3").Value = 2258
